# Update Q3-2025 (row 29) metrics for "bibi" recurrence data
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C29").Value = 24
$ws.Range("D29").Value = 6
$ws.Range("E29").Value = 18
$ws.Range("F29").Value = 1.032702237521515
